$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows("31:31").Insert()

$ws.Range("A31").Value = 2
$ws.Range("B31").Value = 'Comercializadora del Agro de Limarí'
$ws.Range("C31").Value = 'Coquimbo'
$ws.Range("D31").Value = 45036
$ws.Range("E31").Value = 4
$ws.Range("F31").Value = 100112032
$ws.Range("G31").Value = 'Zapallo italiano'
$ws.Range("H31").Value = 'Sin especificar'
$ws.Range("I31").Value = 'Primera'
$ws.Range("J31").Value = 500
$ws.Range("K31").Value = 7000
$ws.Range("L31").Value = 7500
$ws.Range("M31").Value = 7250
$ws.Range("N31").Value = '$/caja 60 unidades'
$ws.Range("O31").Value = 'Provincia de Limarí'
$ws.Range("P31").Value = 121
$ws.Range("Q31").Value = 60
$ws.Range("R31").Value = 'Hortaliza'
